# Weekly Fruta/Hortaliza update: insert a new (most recent) week's record
# at the top of the data block (row 33) and push the existing records
# down by one row, so the oldest existing record (old row 68) is kept
# and re-appears as the new last row (row 69).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a blank row at 33; this shifts old rows 33..68 down to 34..69,
# automatically carrying all of their values/styles along (including the
# previous last row, which becomes the new row 69), and also grows the
# sheet's used range/dimension to A1:R69.
$ws.Rows("33:33").Insert()

# Populate the freshly inserted row 33 with the new week's data. Columns
# A, B, C, E, F, G, H, I, N, O, Q, R are constant across this data block,
# so they carry the same values as every other row; D, J, K, L, M, P are
# the new observations for this row.
$ws.Range("A33").Value = 8
$ws.Range("B33").Value = "Terminal La Palmera de La Serena"
$ws.Range("C33").Value = "Coquimbo"
$ws.Range("D33").Value = 44413
$ws.Range("E33").Value = 4
$ws.Range("F33").Value = 100112044
$ws.Range("G33").Value = "Perejil"
$ws.Range("H33").Value = "Sin especificar"
$ws.Range("I33").Value = "Primera"
$ws.Range("J33").Value = 3360
$ws.Range("K33").Value = 2000
$ws.Range("L33").Value = 2500
$ws.Range("M33").Value = 2250
$ws.Range("N33").Value = "$/atado 1 a 1,5 kilos"
$ws.Range("O33").Value = "Provincia del Elquí"
$ws.Range("P33").Value = 1500
$ws.Range("Q33").Value = 1.5
$ws.Range("R33").Value = "Hortaliza"
